$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.565.04"
$ws.Range("E2").Value = "  -5.59%  "
$ws.Range("D3").Value = "3.360.64"
$ws.Range("E3").Value = "  -6.76%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "182.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.601"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -4.17%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "3.352.42"
$ws.Range("E9").Value = "  -6.61%  "
$ws.Range("E10").Value = "  -12.12%  "
$ws.Range("E11").Value = "  -7.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.06"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.67%  "
$ws.Range("E13").Value = "  -10.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.73"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -9.86%  "
$ws.Range("D15").Value = "3.900.74"
$ws.Range("E15").Value = "  -6.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "608.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -11.20%  "
$ws.Range("D17").Value = "66.526.43"
$ws.Range("E17").Value = "  -5.73%  "
$ws.Range("D18").Value = "3.377.64"
$ws.Range("E18").Value = "  -6.46%  "
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -8.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.912"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -8.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "16.83"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -12.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.08"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  -8.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -10.67%  "
$ws.Range("E32").Value = "  -13.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.84%  "
$ws.Range("B35").Value = "Maker"
$ws.Range("C35").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D35").Value = "3.828.26"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.78"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "534.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.95%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.82"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +38.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.41%  "
$ws.Range("D41").Value = "0.0₃0727"
$ws.Range("E41").Value = "  -14.32%  "
$ws.Range("E42").Value = "  -9.51%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.353"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.81%  "
$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.128"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.62%  "
$ws.Range("E46").Value = "  -10.83%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -12.46%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -9.25%  "
$ws.Range("E49").Value = "  -6.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.21%  "
